$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header for product_collection custom line item property
$ws.Range("E1").Value = "product_collection"

# Populate product_collection values.
# Rows 2 and 3 mirror the product_title column (A) directly.
$ws.Range("E2").Value = $ws.Range("A2").Value2
$ws.Range("E3").Value = $ws.Range("A3").Value2

# Rows 4 and 5 are the "alternate" pairing - swapped relative to column A
# (Modern Muse - 3 Item / Modern Muse - 5 Item are cross-linked as each
# other's alternate product/collection).
$ws.Range("E4").Value = $ws.Range("A5").Value2
$ws.Range("E5").Value = $ws.Range("A4").Value2

# Update the active selection to reflect where the editor left off.
$ws.Range("E10").Select()
